$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "SIM": complete the existing order-6 row (date + "Device" mode) and
# append a new order-7 row underneath it.
# ---------------------------------------------------------------------------
$sim = $wb.Worksheets.Item("SIM")

$sim.Range("D7").Value = 44289
$sim.Range("D7").NumberFormat = "d-mmm-yy"
$sim.Range("F7").Value = "Device"

$sim.Range("A8").Value = 7
$sim.Range("B8").Value = 27000641
$sim.Range("C8").Value = "Connectivity-SIM & CAF"
$sim.Range("D8").Value = 44313
$sim.Range("D8").NumberFormat = "d-mmm-yy"
$sim.Range("E8").Value = 5000
$sim.Range("F8").Value = "Device"

# ---------------------------------------------------------------------------
# Sheet "Phone": append two more "Net from Axis" transfer rows to the right
# hand Reference Id / Date / Amount / Mode table.
# ---------------------------------------------------------------------------
$phone = $wb.Worksheets.Item("Phone")

$phone.Range("J9").Value = 7
$phone.Range("L9").Value = 44338
$phone.Range("L9").NumberFormat = '[$-409]d\-mmm\-yyyy;@'
$phone.Range("M9").Value = 75000
$phone.Range("N9").Value = "Net from Axis"

$phone.Range("J10").Value = 8
$phone.Range("L10").Value = 44342
$phone.Range("L10").NumberFormat = '[$-409]d\-mmm\-yyyy;@'
$phone.Range("M10").Value = 71773
$phone.Range("N10").Value = "Net from Axis"

$null = $phone.Outline.ShowLevels(0, 0)

$null = $phone.Range("E8").Select()

# Re-select on "SIM" last so it remains the active/front-most tab (matches
# the workbook's original tabSelected state).
$null = $sim.Range("E9").Select()
